$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The first two evaluation columns (Öd1/Öd2) are being dropped - the
# remaining table now only tracks Quiz / Vize / Fin, so remove those
# two columns entirely (B:C). This shifts the old D:G columns left to
# become the new B:E columns.
$ws.Range("B1:C8").EntireColumn.Delete()

# Fix up the point-distribution row (row 2) - Quiz/Vize/Fin max points.
$ws.Range("B2").Value = 30
$ws.Range("C2").Value = 30
$ws.Range("D2").Value = 40

# Recompute the TOPLAM column (now column E) as the sum of the three
# remaining score columns (B:D) for each student/outcome row.
$ws.Range("E4").Value = 1
$ws.Range("E5").Value = 2
$ws.Range("E6").Value = 2
$ws.Range("E7").Value = 2
$ws.Range("E8").Value = 3
